{"js": "// Replace each three-digit-by-one-digit multiplication equation in the\n// answers table with its updated value, matching the diff exactly.\n// Each old equation string is unique in the document, so a direct\n// search + Replace on the exact text is safe and precise.\nconst pairs = [\n  [\"312\u00d72=624\", \"896\u00d72=1792\"],\n  [\"356\u00d74=1424\", \"620\u00d72=1240\"],\n  [\"394\u00d76=2364\", \"992\u00d72=1984\"],\n  [\"409\u00d74=1636\", \"684\u00d77=4788\"],\n  [\"426\u00d72=852\", \"219\u00d79=1971\"],\n  [\"490\u00d78=3920\", \"128\u00d72=256\"],\n  [\"500\u00d79=4500\", \"473\u00d77=3311\"],\n  [\"684\u00d73=2052\", \"592\u00d74=2368\"],\n  [\"476\u00d74=1904\", \"492\u00d72=984\"],\n  [\"300\u00d73=900\", \"881\u00d79=7929\"],\n  [\"274\u00d77=1918\", \"108\u00d74=432\"],\n  [\"187\u00d75=935\", \"987\u00d74=3948\"],\n  [\"930\u00d74=3720\", \"494\u00d77=3458\"],\n  [\"905\u00d79=8145\", \"279\u00d77=1953\"],\n  [\"954\u00d76=5724\", \"595\u00d74=2380\"],\n  [\"797\u00d72=1594\", \"391\u00d74=1564\"],\n  [\"513\u00d75=2565\", \"442\u00d74=1768\"],\n  [\"400\u00d73=1200\", \"170\u00d79=1530\"],\n  [\"586\u00d78=4688\", \"945\u00d78=7560\"],\n  [\"941\u00d73=2823\", \"460\u00d74=1840\"],\n  [\"536\u00d74=2144\", \"297\u00d77=2079\"],\n  [\"392\u00d77=2744\", \"395\u00d74=1580\"],\n  [\"741\u00d73=2223\", \"159\u00d73=477\"],\n  [\"842\u00d77=5894\", \"733\u00d79=6597\"],\n  [\"346\u00d77=2422\", \"982\u00d74=3928\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text to replace: \" + oldText);\n  }\n\n  results.items[0].insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n\n", "ps1": "# Update each \"three-digit number x one-digit number\" equation in the\n# answers table to its new value. Every \"old\" equation text is unique\n# within the document, so Find/Execute on the literal text is precise\n# and will not touch any other content (e.g. the date heading).\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @('312\u00d72=624', '896\u00d72=1792'),\n    @('356\u00d74=1424', '620\u00d72=1240'),\n    @('394\u00d76=2364', '992\u00d72=1984'),\n    @('409\u00d74=1636', '684\u00d77=4788'),\n    @('426\u00d72=852', '219\u00d79=1971'),\n    @('490\u00d78=3920', '128\u00d72=256'),\n    @('500\u00d79=4500', '473\u00d77=3311'),\n    @('684\u00d73=2052', '592\u00d74=2368'),\n    @('476\u00d74=1904', '492\u00d72=984'),\n    @('300\u00d73=900', '881\u00d79=7929'),\n    @('274\u00d77=1918', '108\u00d74=432'),\n    @('187\u00d75=935', '987\u00d74=3948'),\n    @('930\u00d74=3720', '494\u00d77=3458'),\n    @('905\u00d79=8145', '279\u00d77=1953'),\n    @('954\u00d76=5724', '595\u00d74=2380'),\n    @('797\u00d72=1594', '391\u00d74=1564'),\n    @('513\u00d75=2565', '442\u00d74=1768'),\n    @('400\u00d73=1200', '170\u00d79=1530'),\n    @('586\u00d78=4688', '945\u00d78=7560'),\n    @('941\u00d73=2823', '460\u00d74=1840'),\n    @('536\u00d74=2144', '297\u00d77=2079'),\n    @('392\u00d77=2744', '395\u00d74=1580'),\n    @('741\u00d73=2223', '159\u00d73=477'),\n    @('842\u00d77=5894', '733\u00d79=6597'),\n    @('346\u00d77=2422', '982\u00d74=3928'),\n    )\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Could not find text to replace: $oldText\"\n    }\n}\n\n"}
